$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,39
$row2[0,0] = 24.07513281277198
$row2[0,1] = 3.799269055677978
$row2[0,2] = 20.275863757094
$row2[0,3] = 52.7481845603303
$row2[0,4] = 52.28133421498961
$row2[0,5] = 21.4525929128096
$row2[0,6] = 18.4756598134299
$row2[0,7] = 1.775636829763065
$row2[0,8] = 0.9736277482982219
$row2[0,9] = 6.5490997487709
$row2[0,10] = 3.17252263552265
$row2[0,11] = 1.726716817859081
$row2[0,12] = 11.06630808873788
$row2[0,13] = 98.99958341210348
$row2[0,14] = 19.52052953313373
$row2[0,15] = 79.47905387896975
$row2[0,16] = 71.61840760179791
$row2[0,17] = 72.12293448120039
$row2[0,18] = 85.7478292417571
$row2[0,19] = 70.9460325627859
$row2[0,20] = 7.034480936064599
$row2[0,21] = 3.539943673163451
$row2[0,22] = 18.96347301816437
$row2[0,23] = 16.64702067334411
$row2[0,24] = 8.234340299792436
$row2[0,25] = 29.44491525312177
$row2[0,26] = 122.8252504052022
$row2[0,27] = 30.52628880454227
$row2[0,28] = 92.29896160065991
$row2[0,29] = 68.13955103761599
$row2[0,30] = 70.79985279798697
$row2[0,31] = 86.08387796585318
$row2[0,32] = 69.43069729599193
$row2[0,33] = 8.881984683447522
$row2[0,34] = 5.007823244097714
$row2[0,35] = 25.03099732277546
$row2[0,36] = 17.98555016221655
$row2[0,37] = 9.132307093088723
$row2[0,38] = 31.25155022418477
$ws.Range("B2:AN2").Value = $row2

$row3 = New-Object 'object[,]' 1,39
$row3[0,0] = 14.07693467630521
$row3[0,1] = 3.140697953161452
$row3[0,2] = 10.93623672314376
$row3[0,3] = 43.34970928390789
$row3[0,4] = 41.12535439892765
$row3[0,5] = 11.92165345075503
$row3[0,6] = 9.608309861400002
$row3[0,7] = 0.9242797126435862
$row3[0,8] = 0.5856052531994604
$row3[0,9] = 4.555358248659908
$row3[0,10] = 2.438476536814568
$row3[0,11] = 1.403645570738901
$row3[0,12] = 9.372275009640541
$row3[0,13] = 76.18205305401318
$row3[0,14] = 20.03490384510583
$row3[0,15] = 56.14714920890734
$row3[0,16] = 65.67357260336011
$row3[0,17] = 63.3409090916464
$row3[0,18] = 61.16687639398653
$row3[0,19] = 46.30127959389502
$row3[0,20] = 5.829050887919577
$row3[0,21] = 2.760658438636084
$row3[0,22] = 17.16546349261412
$row3[0,23] = 17.43541861792866
$row3[0,24] = 8.171686409444012
$row3[0,25] = 31.24002819566737
$row3[0,26] = 97.81811959896311
$row3[0,27] = 28.95409644108173
$row3[0,28] = 68.86402315788136
$row3[0,29] = 65.30251144221495
$row3[0,30] = 65.56629687315132
$row3[0,31] = 65.79117613859749
$row3[0,32] = 49.52266359897934
$row3[0,33] = 6.708004746024338
$row3[0,34] = 3.584143593051416
$row3[0,35] = 19.92525358086478
$row3[0,36] = 18.15556096481395
$row3[0,37] = 8.682109977130404
$row3[0,38] = 31.84960950287157
$ws.Range("B3:AN3").Value = $row3

$row4 = New-Object 'object[,]' 1,39
$row4[0,0] = 22.75815997716055
$row4[0,1] = 4.268021780432877
$row4[0,2] = 18.49013819672767
$row4[0,3] = 53.98322844439765
$row4[0,4] = 50.51176737270973
$row4[0,5] = 20.07939748605907
$row4[0,6] = 16.58559655629098
$row4[0,7] = 2.399180835380752
$row4[0,8] = 1.428700891235248
$row4[0,9] = 9.720047529456998
$row4[0,10] = 3.78164516776528
$row4[0,11] = 1.935356446213969
$row4[0,12] = 11.80065608196892
$row4[0,13] = 98.8990431139448
$row4[0,14] = 19.98930647061056
$row4[0,15] = 78.90973664333424
$row4[0,16] = 71.79870072650439
$row4[0,17] = 72.23971731581183
$row4[0,18] = 85.79407392181928
$row4[0,19] = 69.69741256006049
$row4[0,20] = 6.57367329874972
$row4[0,21] = 3.279842276285141
$row4[0,22] = 18.00385351754227
$row4[0,23] = 16.03870469448468
$row4[0,24] = 7.408123614861689
$row4[0,25] = 29.57120512564511
$row4[0,26] = 118.4716640820771
$row4[0,27] = 28.22192909471696
$row4[0,28] = 90.24973498736018
$row4[0,29] = 69.9864686633582
$row4[0,30] = 72.30759175283009
$row4[0,31] = 86.47563371513718
$row4[0,32] = 69.63355394537818
$row4[0,33] = 8.821969820806668
$row4[0,34] = 5.316936155794318
$row4[0,35] = 24.50269731567888
$row4[0,36] = 17.42991194155889
$row4[0,37] = 8.754167657383242
$row4[0,38] = 30.27629459780395
$ws.Range("B4:AN4").Value = $row4

$row5 = New-Object 'object[,]' 1,39
$row5[0,0] = 17.31328867580875
$row5[0,1] = 3.446624814503255
$row5[0,2] = 13.86666386130549
$row5[0,3] = 49.81878157280467
$row5[0,4] = 45.87855325658418
$row5[0,5] = 14.90779404093295
$row5[0,6] = 12.14679171163789
$row5[0,7] = 1.665514148863828
$row5[0,8] = 1.003786709795077
$row5[0,9] = 6.635274475816622
$row5[0,10] = 2.901381095792405
$row5[0,11] = 1.60396110725102
$row5[0,12] = 9.762550937582409
$row5[0,13] = 91.88269225245652
$row5[0,14] = 17.99808077023203
$row5[0,15] = 73.88461148222449
$row5[0,16] = 71.88005486171394
$row5[0,17] = 72.28136389808543
$row5[0,18] = 78.74822390873091
$row5[0,19] = 65.26194228648809
$row5[0,20] = 6.197778368548109
$row5[0,21] = 3.145444536821154
$row5[0,22] = 16.80860272956033
$row5[0,23] = 15.22341725587524
$row5[0,24] = 7.482918636869745
$row5[0,25] = 28.30862929376817
$row5[0,26] = 103.738027053278
$row5[0,27] = 27.31163912968989
$row5[0,28] = 76.42638792358812
$row5[0,29] = 68.83280628590063
$row5[0,30] = 71.7755949420411
$row5[0,31] = 69.8356768671453
$row5[0,32] = 54.66416789627829
$row5[0,33] = 6.976985106788038
$row5[0,34] = 3.779921217811707
$row5[0,35] = 18.12583320447203
$row5[0,36] = 14.93615805318249
$row5[0,37] = 7.219197042366157
$row5[0,38] = 29.92125870049438
$ws.Range("B5:AN5").Value = $row5

$row6 = New-Object 'object[,]' 1,39
$row6[0,0] = 17.69065414799485
$row6[0,1] = 3.291036990809096
$row6[0,2] = 14.39961715718575
$row6[0,3] = 48.52323533874714
$row6[0,4] = 45.67772324610054
$row6[0,5] = 14.93642499275016
$row6[0,6] = 12.34963915168303
$row6[0,7] = 1.200261313217506
$row6[0,8] = 0.6024001673942236
$row6[0,9] = 4.060414537010583
$row6[0,10] = 2.029669268458242
$row6[0,11] = 0.8823908303144843
$row6[0,12] = 5.234563615903987
$row6[0,13] = 90.37661728890112
$row6[0,14] = 18.47419953629263
$row6[0,15] = 71.90241775260847
$row6[0,16] = 67.37605265102367
$row6[0,17] = 66.04650603646262
$row6[0,18] = 75.99150472750929
$row6[0,19] = 62.58081969756511
$row6[0,20] = 5.551726944994943
$row6[0,21] = 2.870929708728724
$row6[0,22] = 14.8959700397898
$row6[0,23] = 15.14389571625846
$row6[0,24] = 7.66847168972606
$row6[0,25] = 31.12920485523972
$row6[0,26] = 111.58025428998
$row6[0,27] = 29.64408553376758
$row6[0,28] = 81.93616875621245
$row6[0,29] = 66.03405973739379
$row6[0,30] = 66.93951039771228
$row6[0,31] = 75.56212288847425
$row6[0,32] = 58.53417192212706
$row6[0,33] = 8.873870133497777
$row6[0,34] = 4.897245727921488
$row6[0,35] = 22.51870377771627
$row6[0,36] = 17.4120453135647
$row6[0,37] = 8.274708080253724
$row6[0,38] = 30.98993570297703
$ws.Range("B6:AN6").Value = $row6

$row7 = New-Object 'object[,]' 1,39
$row7[0,0] = 17.6766387061932
$row7[0,1] = 4.008119246520345
$row7[0,2] = 13.66851945967285
$row7[0,3] = 46.26753591197575
$row7[0,4] = 41.61771021627266
$row7[0,5] = 14.36524049308764
$row7[0,6] = 11.58249572270408
$row7[0,7] = 1.841034093971801
$row7[0,8] = 1.158882571568058
$row7[0,9] = 7.056863267887542
$row7[0,10] = 3.159746132360045
$row7[0,11] = 1.655012861875784
$row7[0,12] = 10.12127729716251
$row7[0,13] = 91.08559641347405
$row7[0,14] = 18.32789253604562
$row7[0,15] = 72.75770387742844
$row7[0,16] = 69.71125507691873
$row7[0,17] = 70.48181057874424
$row7[0,18] = 75.91539359302277
$row7[0,19] = 62.48066640770838
$row7[0,20] = 6.025427094891376
$row7[0,21] = 3.245937041697988
$row7[0,22] = 16.84893762132593
$row7[0,23] = 14.79594793083408
$row7[0,24] = 7.91188526640682
$row7[0,25] = 31.21139366936779
$row7[0,26] = 119.1534976694762
$row7[0,27] = 29.00183025978855
$row7[0,28] = 90.15166740968763
$row7[0,29] = 69.46179749391739
$row7[0,30] = 70.98224918728029
$row7[0,31] = 78.55054870346187
$row7[0,32] = 63.21412786933215
$row7[0,33] = 8.686427617045508
$row7[0,34] = 4.673613508246891
$row7[0,35] = 21.84330065551743
$row7[0,36] = 18.09927417978575
$row7[0,37] = 9.690524946260044
$row7[0,38] = 34.88394688000852
$ws.Range("B7:AN7").Value = $row7

$row8 = New-Object 'object[,]' 1,39
$row8[0,0] = 16.36553044029763
$row8[0,1] = 3.427576027033513
$row8[0,2] = 12.93795441326412
$row8[0,3] = 48.09788337921533
$row8[0,4] = 42.93150580651595
$row8[0,5] = 13.63792621173332
$row8[0,6] = 11.02466942013299
$row8[0,7] = 1.672722094242557
$row8[0,8] = 0.9611445255345529
$row8[0,9] = 6.686710405083309
$row8[0,10] = 2.691158936512944
$row8[0,11] = 1.262597381255829
$row8[0,12] = 8.019033493675451
$row8[0,13] = 86.47940395269076
$row8[0,14] = 17.77633125201909
$row8[0,15] = 68.70307270067167
$row8[0,16] = 70.25227125226647
$row8[0,17] = 69.08643624741396
$row8[0,18] = 71.38924450787508
$row8[0,19] = 58.36083683979436
$row8[0,20] = 5.470855838333485
$row8[0,21] = 2.853127151154708
$row8[0,22] = 15.76585544746328
$row8[0,23] = 13.05548962941689
$row8[0,24] = 6.369578367805889
$row8[0,25] = 27.11358317918782
$row8[0,26] = 113.3396235616286
$row8[0,27] = 30.58447194797705
$row8[0,28] = 82.75515161365159
$row8[0,29] = 64.85898664584644
$row8[0,30] = 66.62553287204915
$row8[0,31] = 71.66354544775817
$row8[0,32] = 55.81355793746374
$row8[0,33] = 8.197239159856704
$row8[0,34] = 4.633309848843711
$row8[0,35] = 22.93158170297162
$row8[0,36] = 15.16182744409013
$row8[0,37] = 7.62552090187094
$row8[0,38] = 27.76292337135188
$ws.Range("B8:AN8").Value = $row8

$row9 = New-Object 'object[,]' 1,39
$row9[0,0] = 16.96481753582974
$row9[0,1] = 3.1776091974884
$row9[0,2] = 13.78720833834134
$row9[0,3] = 45.53759745090429
$row9[0,4] = 42.21743005313715
$row9[0,5] = 13.4613896598395
$row9[0,6] = 11.24247467874589
$row9[0,7] = 0.9893092848932264
$row9[0,8] = 0.6005517763154823
$row9[0,9] = 4.503095661929195
$row9[0,10] = 2.537890999151071
$row9[0,11] = 1.452136872352977
$row9[0,12] = 8.975423112845307
$row9[0,13] = 85.18571043899587
$row9[0,14] = 18.86675201319139
$row9[0,15] = 66.31895842580448
$row9[0,16] = 65.82178011685285
$row9[0,17] = 66.9922929044394
$row9[0,18] = 69.08425165729618
$row9[0,19] = 56.10481233262706
$row9[0,20] = 5.591443523431761
$row9[0,21] = 3.186530913494789
$row9[0,22] = 17.12029421405072
$row9[0,23] = 13.62879098631295
$row9[0,24] = 7.135166846111168
$row9[0,25] = 29.44448686639759
$row9[0,26] = 116.8066656059019
$row9[0,27] = 32.13606828910684
$row9[0,28] = 84.67059731679505
$row9[0,29] = 66.42450132842434
$row9[0,30] = 69.21792205654432
$row9[0,31] = 72.99337281411715
$row9[0,32] = 57.00951854649204
$row9[0,33] = 7.722408775931791
$row9[0,34] = 3.616513595194894
$row9[0,35] = 16.39769276759867
$row9[0,36] = 14.2690513855833
$row9[0,37] = 7.009110157831941
$row9[0,38] = 28.92465362135869
$ws.Range("B9:AN9").Value = $row9

$row10 = New-Object 'object[,]' 1,39
$row10[0,0] = 14.69104865697746
$row10[0,1] = 3.345570840335869
$row10[0,2] = 11.3454778166416
$row10[0,3] = 43.15980036677967
$row10[0,4] = 36.59426292538839
$row10[0,5] = 11.61274819390905
$row10[0,6] = 9.105708333795565
$row10[0,7] = 1.124820729882877
$row10[0,8] = 0.6388839708073031
$row10[0,9] = 4.198027037554691
$row10[0,10] = 2.376611449552886
$row10[0,11] = 1.157688248530637
$row10[0,12] = 7.007346291585359
$row10[0,13] = 83.31841821607721
$row10[0,14] = 18.37586528914189
$row10[0,15] = 64.94255292693531
$row10[0,16] = 68.22683507895383
$row10[0,17] = 68.20099892285242
$row10[0,18] = 68.38642129966222
$row10[0,19] = 55.17287341130149
$row10[0,20] = 4.173845583938286
$row10[0,21] = 2.058592226850599
$row10[0,22] = 12.23433408509346
$row10[0,23] = 12.73518897048659
$row10[0,24] = 6.156406243903709
$row10[0,25] = 29.09349128590032
$row10[0,26] = 112.122782082962
$row10[0,27] = 33.36636109440304
$row10[0,28] = 78.75642098855899
$row10[0,29] = 62.30788874421224
$row10[0,30] = 64.77121404253775
$row10[0,31] = 72.28504600379124
$row10[0,32] = 54.06717515891595
$row10[0,33] = 6.975656709461317
$row10[0,34] = 3.41784078247727
$row10[0,35] = 17.73361697697374
$row10[0,36] = 17.45453363696333
$row10[0,37] = 8.110785324022
$row10[0,38] = 27.89717143185874
$ws.Range("B10:AN10").Value = $row10

$row11 = New-Object 'object[,]' 1,39
$row11[0,0] = 15.3281308607605
$row11[0,1] = 2.807465580446539
$row11[0,2] = 12.52066528031396
$row11[0,3] = 48.74215340293917
$row11[0,4] = 43.8354963942834
$row11[0,5] = 12.34291884975213
$row11[0,6] = 10.32918925894096
$row11[0,7] = 1.514947848148667
$row11[0,8] = 0.858964402043215
$row11[0,9] = 4.944229359648741
$row11[0,10] = 2.391933078123962
$row11[0,11] = 1.35302028114197
$row11[0,12] = 8.371861198521225
$row11[0,13] = 83.58962351126009
$row11[0,14] = 19.85138234216787
$row11[0,15] = 63.73824116909222
$row11[0,16] = 65.36450100877227
$row11[0,17] = 64.82078866231775
$row11[0,18] = 66.08951352994464
$row11[0,19] = 52.20052691194528
$row11[0,20] = 5.809334796126017
$row11[0,21] = 2.905102417303563
$row11[0,22] = 16.55202897402063
$row11[0,23] = 13.5319270623441
$row11[0,24] = 6.24692712540474
$row11[0,25] = 25.81042737637703
$row11[0,26] = 115.3715248015184
$row11[0,27] = 30.55666794909486
$row11[0,28] = 84.81485685242353
$row11[0,29] = 67.03631475329045
$row11[0,30] = 68.27156380639212
$row11[0,31] = 76.0153800440454
$row11[0,32] = 59.5207946738859
$row11[0,33] = 8.219766510622115
$row11[0,34] = 4.308228262097574
$row11[0,35] = 21.45081535483823
$row11[0,36] = 16.75284010540112
$row11[0,37] = 7.827447192101293
$row11[0,38] = 28.66780899446062
$ws.Range("B11:AN11").Value = $row11

$row12 = New-Object 'object[,]' 1,39
$row12[0,0] = 13.97923663452209
$row12[0,1] = 3.812930515452838
$row12[0,2] = 10.16630611906925
$row12[0,3] = 40.7560077103464
$row12[0,4] = 37.65165392896905
$row12[0,5] = 11.60039252544207
$row12[0,6] = 8.613005845597547
$row12[0,7] = 1.535536293102836
$row12[0,8] = 0.8130339479959224
$row12[0,9] = 5.766456558543575
$row12[0,10] = 2.579232603686354
$row12[0,11] = 1.186476347592796
$row12[0,12] = 7.873809854846473
$row12[0,13] = 71.01561559702579
$row12[0,14] = 17.64268380986872
$row12[0,15] = 53.37293178715705
$row12[0,16] = 63.08925492404941
$row12[0,17] = 62.44458438402722
$row12[0,18] = 56.54960456911451
$row12[0,19] = 44.11823454575636
$row12[0,20] = 4.514667215464466
$row12[0,21] = 2.145324180456512
$row12[0,22] = 12.11282430390409
$row12[0,23] = 11.15206479396653
$row12[0,24] = 5.347910592232108
$row12[0,25] = 23.44584627044557
$row12[0,26] = 95.49673473966499
$row12[0,27] = 28.42179308122914
$row12[0,28] = 67.07494165843585
$row12[0,29] = 64.1510175570194
$row12[0,30] = 64.85659945167386
$row12[0,31] = 60.09391484240266
$row12[0,32] = 44.26028821122707
$row12[0,33] = 6.304475189903406
$row12[0,34] = 3.217091735561892
$row12[0,35] = 17.97674630636103
$row12[0,36] = 15.78268431045799
$row12[0,37] = 7.323645424066192
$row12[0,38] = 29.49595163120384
$ws.Range("B12:AN12").Value = $row12

$row13 = New-Object 'object[,]' 1,39
$row13[0,0] = 14.85345876852156
$row13[0,1] = 3.656159914373765
$row13[0,2] = 11.19729885414779
$row13[0,3] = 44.11620962312067
$row13[0,4] = 40.56243175918797
$row13[0,5] = 12.37558413616537
$row13[0,6] = 9.512433759448006
$row13[0,7] = 1.548640622645488
$row13[0,8] = 0.8141351212492267
$row13[0,9] = 5.530115316003487
$row13[0,10] = 2.649762413433432
$row13[0,11] = 1.119983967892213
$row13[0,12] = 6.75280524710749
$row13[0,13] = 78.13218262023491
$row13[0,14] = 19.29515796486013
$row13[0,15] = 58.83702465537479
$row13[0,16] = 66.23750328625175
$row13[0,17] = 66.7382938998781
$row13[0,18] = 64.1419707723925
$row13[0,19] = 49.88286350076933
$row13[0,20] = 5.103698277584169
$row13[0,21] = 2.339441432866868
$row13[0,22] = 13.60468791740218
$row13[0,23] = 14.4148436352268
$row13[0,24] = 6.78669397281808
$row13[0,25] = 28.1520702116596
$row13[0,26] = 106.9142713501164
$row13[0,27] = 31.26098395435612
$row13[0,28] = 75.65328739576024
$row13[0,29] = 62.79706020851042
$row13[0,30] = 65.1099482166578
$row13[0,31] = 69.44251793239852
$row13[0,32] = 52.72918267818524
$row13[0,33] = 6.906435812078462
$row13[0,34] = 3.132099393244749
$row13[0,35] = 16.60899057570892
$row13[0,36] = 15.15012473035076
$row13[0,37] = 6.962303847441522
$row13[0,38] = 29.67424835062988
$ws.Range("B13:AN13").Value = $row13

$row14 = New-Object 'object[,]' 1,39
$row14[0,0] = 14.64555327868944
$row14[0,1] = 3.040648297476291
$row14[0,2] = 11.60490498121315
$row14[0,3] = 49.04139873292924
$row14[0,4] = 43.91861039009556
$row14[0,5] = 11.49181015260686
$row14[0,6] = 9.423327625324927
$row14[0,7] = 1.4829096488668
$row14[0,8] = 0.5689102415298609
$row14[0,9] = 4.009249883432138
$row14[0,10] = 2.6723886811518
$row14[0,11] = 1.317608767042793
$row14[0,12] = 9.136181158132421
$row14[0,13] = 81.60507286262894
$row14[0,14] = 17.50231452912452
$row14[0,15] = 64.10275833350443
$row14[0,16] = 81.19524844424498
$row14[0,17] = 80.44985770118041
$row14[0,18] = 63.84952299004522
$row14[0,19] = 51.57156768142003
$row14[0,20] = 3.558480335384085
$row14[0,21] = 1.649589213809967
$row14[0,22] = 11.04955721990763
$row14[0,23] = 13.3740847878556
$row14[0,24] = 6.864874256789134
$row14[0,25] = 36.42220136471755
$row14[0,26] = 103.0691820818958
$row14[0,27] = 27.72194873773205
$row14[0,28] = 75.34723334416377
$row14[0,29] = 77.28438484741231
$row14[0,30] = 79.05031600778267
$row14[0,31] = 63.89805917286302
$row14[0,32] = 50.36132558369771
$row14[0,33] = 5.657304976056472
$row14[0,34] = 3.073684147200285
$row14[0,35] = 17.24276396888647
$row14[0,36] = 14.14998484865068
$row14[0,37] = 7.408592692752195
$row14[0,38] = 33.87376505326314
$ws.Range("B14:AN14").Value = $row14

$row15 = New-Object 'object[,]' 1,39
$row15[0,0] = 18.88279809212263
$row15[0,1] = 3.910874836704027
$row15[0,2] = 14.9719232554186
$row15[0,3] = 49.95944009372553
$row15[0,4] = 46.18617325593083
$row15[0,5] = 15.99573163917673
$row15[0,6] = 12.90998678313631
$row15[0,7] = 1.556502204811441
$row15[0,8] = 0.954745087033036
$row15[0,9] = 5.85432685979562
$row15[0,10] = 2.743374517423797
$row15[0,11] = 1.420300797303592
$row15[0,12] = 9.537826722397339
$row15[0,13] = 93.73749248589787
$row15[0,14] = 18.90583269173788
$row15[0,15] = 74.83165979415998
$row15[0,16] = 72.27407356419398
$row15[0,17] = 73.45397586389647
$row15[0,18] = 78.01004265469429
$row15[0,19] = 64.62430652348286
$row15[0,20] = 5.662165802182868
$row15[0,21] = 3.447789514642229
$row15[0,22] = 17.84346896683171
$row15[0,23] = 14.30581982226358
$row15[0,24] = 7.996509919765154
$row15[0,25] = 31.38641387907116
$row15[0,26] = 121.380398610189
$row15[0,27] = 30.95732266051115
$row15[0,28] = 90.42307594967787
$row15[0,29] = 69.23092412206947
$row15[0,30] = 72.70217562564625
$row15[0,31] = 78.67156757567868
$row15[0,32] = 62.70640436299025
$row15[0,33] = 7.093905439252747
$row15[0,34] = 3.826223993521194
$row15[0,35] = 16.25445947655175
$row15[0,36] = 14.48303187504026
$row15[0,37] = 6.874719302221772
$row15[0,38] = 28.91036175690352
$ws.Range("B15:AN15").Value = $row15

$row16 = New-Object 'object[,]' 1,39
$row16[0,0] = 14.81238114364698
$row16[0,1] = 3.703413419636287
$row16[0,2] = 11.10896772401069
$row16[0,3] = 43.0627211499372
$row16[0,4] = 40.36690994495248
$row16[0,5] = 12.23209862216916
$row16[0,6] = 9.530319098926867
$row16[0,7] = 1.415384056806869
$row16[0,8] = 0.7317066836512424
$row16[0,9] = 4.961736391570249
$row16[0,10] = 2.604972513711148
$row16[0,11] = 1.185323577009718
$row16[0,12] = 7.436762465434569
$row16[0,13] = 81.49657853124839
$row16[0,14] = 19.65824047612619
$row16[0,15] = 61.83833805512221
$row16[0,16] = 65.10413067387215
$row16[0,17] = 64.55552528016173
$row16[0,18] = 66.8220542440252
$row16[0,19] = 52.84068537541955
$row16[0,20] = 4.078290103952602
$row16[0,21] = 2.042953475971775
$row16[0,22] = 12.30589590583069
$row16[0,23] = 13.29489520946396
$row16[0,24] = 6.2908116780883
$row16[0,25] = 27.24254366266403
$row16[0,26] = 109.887855113836
$row16[0,27] = 30.23798948832988
$row16[0,28] = 79.64986562550615
$row16[0,29] = 66.45079243917263
$row16[0,30] = 66.7395038582837
$row16[0,31] = 70.12744795049038
$row16[0,32] = 53.7311430837851
$row16[0,33] = 6.195946543101891
$row16[0,34] = 3.261245510589537
$row16[0,35] = 16.41122960020231
$row16[0,36] = 15.54606962889978
$row16[0,37] = 6.66011110046747
$row16[0,38] = 24.80964694110689
$ws.Range("B16:AN16").Value = $row16

$row17 = New-Object 'object[,]' 1,39
$row17[0,0] = 14.67429325044327
$row17[0,1] = 3.634477527125791
$row17[0,2] = 11.03981572331748
$row17[0,3] = 42.88151965463096
$row17[0,4] = 39.6540339530249
$row17[0,5] = 11.73267012750073
$row17[0,6] = 9.132338464237499
$row17[0,7] = 1.376982183487157
$row17[0,8] = 0.7543443988783992
$row17[0,9] = 5.420325935556107
$row17[0,10] = 2.946124901102464
$row17[0,11] = 1.521332225878391
$row17[0,12] = 9.661371786834023
$row17[0,13] = 79.16302659632146
$row17[0,14] = 20.39980137870329
$row17[0,15] = 58.76322521761816
$row17[0,16] = 65.73616963169314
$row17[0,17] = 66.21126302242968
$row17[0,18] = 61.92062304225777
$row17[0,19] = 47.36068669535118
$row17[0,20] = 4.916266460913808
$row17[0,21] = 2.329225823030574
$row17[0,22] = 13.48428208238882
$row17[0,23] = 15.95875498810377
$row17[0,24] = 7.711178323391443
$row17[0,25] = 27.62952233522175
$row17[0,26] = 117.3376318987472
$row17[0,27] = 35.4170047483347
$row17[0,28] = 81.9206271504125
$row17[0,29] = 65.19962198537381
$row17[0,30] = 66.26963351264374
$row17[0,31] = 69.76896006438577
$row17[0,32] = 51.52997947316225
$row17[0,33] = 6.845423091767547
$row17[0,34] = 3.359001937678561
$row17[0,35] = 17.64933852807768
$row17[0,36] = 17.9556642134608
$row17[0,37] = 7.950247589956167
$row17[0,38] = 27.05795912088067
$ws.Range("B17:AN17").Value = $row17

$row18 = New-Object 'object[,]' 1,39
$row18[0,0] = 13.87168136314738
$row18[0,1] = 3.739463915630331
$row18[0,2] = 10.13221744751705
$row18[0,3] = 40.43866046058734
$row18[0,4] = 37.44910468836319
$row18[0,5] = 11.41170923243052
$row18[0,6] = 8.568780964629765
$row18[0,7] = 1.579976645245478
$row18[0,8] = 0.7872674200081322
$row18[0,9] = 4.967797623805043
$row18[0,10] = 2.797048206748694
$row18[0,11] = 1.168430441389049
$row18[0,12] = 6.7336751272746
$row18[0,13] = 68.56018284690543
$row18[0,14] = 16.78804206887002
$row18[0,15] = 51.77214077803541
$row18[0,16] = 64.92407492201059
$row18[0,17] = 65.39255850312774
$row18[0,18] = 56.58705636444088
$row18[0,19] = 43.74543878197437
$row18[0,20] = 3.759714581620244
$row18[0,21] = 1.579572789200789
$row18[0,22] = 10.14282693545266
$row18[0,23] = 13.73005397094263
$row18[0,24] = 6.70370725956381
$row18[0,25] = 28.20783595140335
$row18[0,26] = 99.29746690811595
$row18[0,27] = 29.21659503880936
$row18[0,28] = 70.08087186930659
$row18[0,29] = 64.48546547536027
$row18[0,30] = 66.03039784574278
$row18[0,31] = 65.42562252815215
$row18[0,32] = 48.50664170826022
$row18[0,33] = 5.64815265421075
$row18[0,34] = 2.883437717229581
$row18[0,35] = 15.52515889465533
$row18[0,36] = 16.50062918493531
$row18[0,37] = 7.233479440646454
$row18[0,38] = 27.73904479200298
$ws.Range("B18:AN18").Value = $row18

$row19 = New-Object 'object[,]' 1,39
$row19[0,0] = 12.59960084945975
$row19[0,1] = 3.319463388725613
$row19[0,2] = 9.280137460734132
$row19[0,3] = 40.22069751532179
$row19[0,4] = 36.12117739920154
$row19[0,5] = 9.982252237163362
$row19[0,6] = 7.466097074790996
$row19[0,7] = 1.630825637491954
$row19[0,8] = 0.7145581229201305
$row19[0,9] = 4.887636285936378
$row19[0,10] = 2.924690556919917
$row19[0,11] = 1.467448093598964
$row19[0,12] = 9.508448347403073
$row19[0,13] = 70.36474346047352
$row19[0,14] = 20.00441627491179
$row19[0,15] = 50.36032718556174
$row19[0,16] = 63.00980639560355
$row19[0,17] = 61.53785168625087
$row19[0,18] = 55.49314934157447
$row19[0,19] = 40.31330386984196
$row19[0,20] = 4.541174715981398
$row19[0,21] = 2.225785436200061
$row19[0,22] = 13.21502387394781
$row19[0,23] = 15.8842538955385
$row19[0,24] = 7.743098601820691
$row19[0,25] = 29.94498528971125
$row19[0,26] = 108.0016365356318
$row19[0,27] = 35.92875050943569
$row19[0,28] = 72.07288602619607
$row19[0,29] = 60.42032283361985
$row19[0,30] = 60.95868848932685
$row19[0,31] = 65.30097774527013
$row19[0,32] = 46.1363627912244
$row19[0,33] = 7.201114269097665
$row19[0,34] = 3.530385403225213
$row19[0,35] = 18.24931292927612
$row19[0,36] = 18.29753078694606
$row19[0,37] = 8.119843523516687
$row19[0,38] = 29.70519951930848
$ws.Range("B19:AN19").Value = $row19

$row20 = New-Object 'object[,]' 1,39
$row20[0,0] = 14.49681562020026
$row20[0,1] = 3.531575940800614
$row20[0,2] = 10.96523967939965
$row20[0,3] = 39.75817072633362
$row20[0,4] = 36.50101596287576
$row20[0,5] = 11.79652747556604
$row20[0,6] = 9.106809331701635
$row20[0,7] = 1.583668200468898
$row20[0,8] = 0.791635939721721
$row20[0,9] = 5.379541951557614
$row20[0,10] = 3.16517434217535
$row20[0,11] = 1.519618213785302
$row20[0,12] = 9.447972116172147
$row20[0,13] = 77.08263156306516
$row20[0,14] = 19.89660857905861
$row20[0,15] = 57.18602298400656
$row20[0,16] = 66.20113764520659
$row20[0,17] = 65.33417636469544
$row20[0,18] = 62.66659603719806
$row20[0,19] = 47.69094525793832
$row20[0,20] = 4.2924359924432
$row20[0,21] = 1.909356750011376
$row20[0,22] = 11.01956100250181
$row20[0,23] = 16.70828581142854
$row20[0,24] = 7.856919228190004
$row20[0,25] = 28.28251925163655
$row20[0,26] = 113.2564719180125
$row20[0,27] = 34.26534282429611
$row20[0,28] = 78.99112909371638
$row20[0,29] = 62.33376718408424
$row20[0,30] = 63.52993571102538
$row20[0,31] = 73.4314559972867
$row20[0,32] = 53.4791604854675
$row20[0,33] = 6.095583497040339
$row20[0,34] = 3.068648100247041
$row20[0,35] = 14.87773571360609
$row20[0,36] = 19.99448611576302
$row20[0,37] = 9.293396606780552
$row20[0,38] = 30.91542540201051
$ws.Range("B20:AN20").Value = $row20

$row21 = New-Object 'object[,]' 1,39
$row21[0,0] = 16.17079315577051
$row21[0,1] = 3.576648449636548
$row21[0,2] = 12.59414470613396
$row21[0,3] = 48.73575132277645
$row21[0,4] = 45.08046445831232
$row21[0,5] = 13.31649663331179
$row21[0,6] = 10.47955445782757
$row21[0,7] = 1.760929301396694
$row21[0,8] = 0.9192201989917174
$row21[0,9] = 6.409033917361985
$row21[0,10] = 2.922098156804878
$row21[0,11] = 1.37213350499055
$row21[0,12] = 8.705066782566947
$row21[0,13] = 81.61110113888429
$row21[0,14] = 19.70946859440713
$row21[0,15] = 61.90163254447717
$row21[0,16] = 68.86089194727502
$row21[0,17] = 68.39785479103131
$row21[0,18] = 65.80382372288585
$row21[0,19] = 51.36200691435849
$row21[0,20] = 5.872336900334623
$row21[0,21] = 2.685348980223544
$row21[0,22] = 14.8826042863519
$row21[0,23] = 16.16642565857562
$row21[0,24] = 7.866273256220625
$row21[0,25] = 28.58423646786831
$row21[0,26] = 116.8418256420334
$row21[0,27] = 32.98065952056196
$row21[0,28] = 83.86116612147144
$row21[0,29] = 69.30838366547434
$row21[0,30] = 68.82894691773539
$row21[0,31] = 73.26353398620289
$row21[0,32] = 54.93070390637272
$row21[0,33] = 8.901036117796673
$row21[0,34] = 4.279175294235978
$row21[0,35] = 19.91919587516946
$row21[0,36] = 18.40195142160213
$row21[0,37] = 8.624637888349769
$row21[0,38] = 31.57452484418775
$ws.Range("B21:AN21").Value = $row21

